# Recalibrated extrapolation values after removing a less-than-USD-5 price
# data point that was treated as noise in the calibration input.
# Only the derived columns (ABSM1_RN, M1_RN, CM2_RN, CMN3_RN, CMN4_RN) in
# rows whose calibration set included that noisy point were recomputed;
# all other cells (labels, TTM, and the *_PH columns) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = 119173.8124178505
$ws.Range("E3").Value = -0.0165108791718419
$ws.Range("F3").Value = 0.1622768345493927
$ws.Range("G3").Value = -0.9397993120033892
$ws.Range("H3").Value = 9.783183741663269

$ws.Range("D5").Value = 120877.9374610504
$ws.Range("E5").Value = -0.01980176926448701
$ws.Range("F5").Value = 0.2021236676005423
$ws.Range("G5").Value = -0.7896683851472568
$ws.Range("H5").Value = 8.701746013044438

$ws.Range("D6").Value = 121332.9788333115
$ws.Range("E6").Value = -0.03141706436792297
$ws.Range("F6").Value = 0.2395987164369215
$ws.Range("G6").Value = -1.227447440299993
$ws.Range("H6").Value = 10.46161784128113

$ws.Range("D7").Value = 122567.858691734
$ws.Range("E7").Value = -0.03322405798612654
$ws.Range("F7").Value = 0.2411684987456979
$ws.Range("G7").Value = -0.6963525108543913
$ws.Range("H7").Value = 5.579453760602306

$ws.Range("D8").Value = 122976.3775938635
$ws.Range("E8").Value = -0.05130078717971131
$ws.Range("F8").Value = 0.2179065660504709
$ws.Range("G8").Value = -0.81845030295635
$ws.Range("H8").Value = 6.389084136372018

$ws.Range("D9").Value = 124609.9680254168
$ws.Range("E9").Value = -0.07968831576119446
$ws.Range("F9").Value = 0.3294018326843595
$ws.Range("G9").Value = -1.419806925590777
$ws.Range("H9").Value = 9.24521707517148

$ws.Range("D10").Value = 126070.0701046711
$ws.Range("E10").Value = -0.114435350839649
$ws.Range("F10").Value = 0.4294800091287003
$ws.Range("G10").Value = -1.90576394231044
$ws.Range("H10").Value = 9.916455759707791

$ws.Range("D11").Value = 128137.9165053551
$ws.Range("E11").Value = -0.1880243208342523
$ws.Range("F11").Value = 0.7478405327394673
$ws.Range("G11").Value = -2.617708984956492
$ws.Range("H11").Value = 12.86694679065372

$ws.Range("D12").Value = 118396.153735129
$ws.Range("E12").Value = -0.03575016654458062
$ws.Range("F12").Value = 0.1319657899133331
$ws.Range("G12").Value = -1.264644170504328
$ws.Range("H12").Value = 11.0022762263902

$ws.Range("D15").Value = 118464.1871531095
$ws.Range("E15").Value = -0.0277054698797307
$ws.Range("F15").Value = 0.1421213205178096
$ws.Range("G15").Value = -0.3546488050891341
$ws.Range("H15").Value = 5.74444265314017

$ws.Range("D17").Value = 118401.252049103
$ws.Range("E17").Value = -0.03024213019375558
$ws.Range("F17").Value = 0.1376752397513106
$ws.Range("G17").Value = -0.5771489722336368
$ws.Range("H17").Value = 6.480858352249357

$ws.Range("D20").Value = 119583.1890569612
$ws.Range("E20").Value = -0.002226389464698468
$ws.Range("F20").Value = 0.1474509123258346
$ws.Range("G20").Value = -0.3607265038845483
$ws.Range("H20").Value = 7.367402260957602
